$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the invoice number in row 2 (A2) from A19252 to A19260
$ws.Range("A2").Value = "A19260"

# Delete rows 3 to 5 (the extra invoice rows no longer present)
$ws.Range("A3:A5").EntireRow.Delete()

# Update the selection to match the target state
$ws.Range("C9").Select()
